$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Row 19 (AGE_ANTH_FUP) content rewrite.
#    The old rich-text algorithm description is replaced by a short variable
#    name, the rule_category / status_detail move to the "direct mapping /
#    identical" pattern used by the other simple rows, and the now-redundant
#    comment cell (I19) is removed entirely.
# ---------------------------------------------------------------------------
$ws.Range("F19").Value = "age_anth_f4"
$ws.Range("G19").Value = "direct_mapping"
$ws.Range("H19").Value = "direct_mapping"
$ws.Range("H19").Style = "Standard"
$ws.Range("I19").ClearContents()
$ws.Range("K19").Value = "identical"

# ---------------------------------------------------------------------------
# 2. Re-colour the algorithm cells that used the plain "explicit Calibri"
#    style (cellXfs index 3) so they use the red-font style (cellXfs index 2)
#    instead -- this is the font-table cleanup that drops the spare font.
# ---------------------------------------------------------------------------
$redCells = @(
    "F3","F4","F5","F6","H6","F7",
    "F9","F10","F12","F13",
    "F20","F21","F22","F23","F24","F25","F26","F27","F28","H28",
    "F31","F32","F33","F34","F35","F36","H36"
)
foreach ($addr in $redCells) {
    $ws.Range($addr).Font.Color = 255
}

# ---------------------------------------------------------------------------
# 3. Scroll / selection state.
# ---------------------------------------------------------------------------
$ws.Range("C19").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.Zoom = 90
